{"js": "const replacements = [\n  [\"84\u00d712=1008\", \"22\u00d729=638\"],\n  [\"65\u00d723=1495\", \"51\u00d754=2754\"],\n  [\"12\u00d759=708\", \"22\u00d771=1562\"],\n  [\"49\u00d781=3969\", \"15\u00d775=1125\"],\n  [\"48\u00d758=2784\", \"49\u00d775=3675\"],\n  [\"29\u00d743=1247\", \"46\u00d779=3634\"],\n  [\"34\u00d746=1564\", \"34\u00d791=3094\"],\n  [\"95\u00d753=5035\", \"57\u00d777=4389\"],\n  [\"53\u00d774=3922\", \"58\u00d759=3422\"],\n  [\"47\u00d762=2914\", \"87\u00d772=6264\"],\n  [\"33\u00d745=1485\", \"36\u00d776=2736\"],\n  [\"30\u00d752=1560\", \"37\u00d716=592\"],\n  [\"91\u00d784=7644\", \"80\u00d779=6320\"],\n  [\"67\u00d766=4422\", \"51\u00d714=714\"],\n  [\"18\u00d789=1602\", \"82\u00d793=7626\"],\n  [\"47\u00d758=2726\", \"86\u00d737=3182\"],\n  [\"87\u00d771=6177\", \"20\u00d737=740\"],\n  [\"65\u00d764=4160\", \"87\u00d763=5481\"],\n  [\"22\u00d796=2112\", \"43\u00d774=3182\"],\n  [\"81\u00d775=6075\", \"62\u00d712=744\"],\n  [\"25\u00d753=1325\", \"36\u00d753=1908\"],\n  [\"36\u00d763=2268\", \"23\u00d795=2185\"],\n  [\"25\u00d748=1200\", \"87\u00d783=7221\"],\n  [\"42\u00d737=1554\", \"34\u00d794=3196\"],\n  [\"16\u00d732=512\", \"88\u00d718=1584\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Not found: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"84\u00d712=1008\", \"22\u00d729=638\"),\n    @(\"65\u00d723=1495\", \"51\u00d754=2754\"),\n    @(\"12\u00d759=708\", \"22\u00d771=1562\"),\n    @(\"49\u00d781=3969\", \"15\u00d775=1125\"),\n    @(\"48\u00d758=2784\", \"49\u00d775=3675\"),\n    @(\"29\u00d743=1247\", \"46\u00d779=3634\"),\n    @(\"34\u00d746=1564\", \"34\u00d791=3094\"),\n    @(\"95\u00d753=5035\", \"57\u00d777=4389\"),\n    @(\"53\u00d774=3922\", \"58\u00d759=3422\"),\n    @(\"47\u00d762=2914\", \"87\u00d772=6264\"),\n    @(\"33\u00d745=1485\", \"36\u00d776=2736\"),\n    @(\"30\u00d752=1560\", \"37\u00d716=592\"),\n    @(\"91\u00d784=7644\", \"80\u00d779=6320\"),\n    @(\"67\u00d766=4422\", \"51\u00d714=714\"),\n    @(\"18\u00d789=1602\", \"82\u00d793=7626\"),\n    @(\"47\u00d758=2726\", \"86\u00d737=3182\"),\n    @(\"87\u00d771=6177\", \"20\u00d737=740\"),\n    @(\"65\u00d764=4160\", \"87\u00d763=5481\"),\n    @(\"22\u00d796=2112\", \"43\u00d774=3182\"),\n    @(\"81\u00d775=6075\", \"62\u00d712=744\"),\n    @(\"25\u00d753=1325\", \"36\u00d753=1908\"),\n    @(\"36\u00d763=2268\", \"23\u00d795=2185\"),\n    @(\"25\u00d748=1200\", \"87\u00d783=7221\"),\n    @(\"42\u00d737=1554\", \"34\u00d794=3196\"),\n    @(\"16\u00d732=512\", \"88\u00d718=1584\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $result = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $result) {\n        throw \"Replace failed for: $old\"\n    }\n}"}
